{"js": "// Update the date line at the top of the document.\nconst body = context.document.body;\nconst dateResults = body.search(\"2025-10-29 Wednesday\", { matchCase: true });\ndateResults.load('items');\nawait context.sync();\nif (dateResults.items.length === 0) {\n  throw new Error('Date text not found: 2025-10-29 Wednesday');\n}\ndateResults.items[0].insertText(\"2025-10-30 Thursday\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Update each multiplication-problem cell in the table by (row, column)\n// position, scoping the search to that single cell so that newly written\n// values never get matched again later. Several target values duplicate\n// other cells' original values (e.g. row 0/col 2 becomes \"62\\u00d731=\"\n// while row 14/col 1 originally held \"62\\u00d731=\" and becomes something\n// else), so a document-wide search/replace would mismatch occurrences.\nconst tables = context.document.body.tables;\ntables.load('items');\nawait context.sync();\nconst table = tables.items[0];\n\nconst cellReplacements = [\n  { row: 0, col: 0, oldText: \"49\u00d764=\", newText: \"72\u00d735=\" },\n  { row: 0, col: 1, oldText: \"23\u00d738=\", newText: \"52\u00d780=\" },\n  { row: 0, col: 2, oldText: \"40\u00d797=\", newText: \"62\u00d731=\" },\n  { row: 0, col: 3, oldText: \"58\u00d786=\", newText: \"84\u00d736=\" },\n  { row: 0, col: 4, oldText: \"75\u00d736=\", newText: \"99\u00d772=\" },\n  { row: 4, col: 0, oldText: \"64\u00d747=\", newText: \"28\u00d738=\" },\n  { row: 4, col: 1, oldText: \"18\u00d714=\", newText: \"19\u00d745=\" },\n  { row: 4, col: 2, oldText: \"89\u00d722=\", newText: \"48\u00d755=\" },\n  { row: 4, col: 3, oldText: \"54\u00d743=\", newText: \"80\u00d778=\" },\n  { row: 4, col: 4, oldText: \"89\u00d745=\", newText: \"89\u00d739=\" },\n  { row: 9, col: 0, oldText: \"31\u00d724=\", newText: \"31\u00d781=\" },\n  { row: 9, col: 1, oldText: \"33\u00d787=\", newText: \"88\u00d799=\" },\n  { row: 9, col: 2, oldText: \"71\u00d781=\", newText: \"27\u00d735=\" },\n  { row: 9, col: 3, oldText: \"36\u00d799=\", newText: \"26\u00d767=\" },\n  { row: 9, col: 4, oldText: \"49\u00d770=\", newText: \"28\u00d719=\" },\n  { row: 14, col: 0, oldText: \"32\u00d720=\", newText: \"92\u00d798=\" },\n  { row: 14, col: 1, oldText: \"62\u00d731=\", newText: \"34\u00d771=\" },\n  { row: 14, col: 2, oldText: \"63\u00d762=\", newText: \"18\u00d797=\" },\n  { row: 14, col: 3, oldText: \"12\u00d790=\", newText: \"80\u00d724=\" },\n  { row: 14, col: 4, oldText: \"25\u00d773=\", newText: \"86\u00d791=\" },\n  { row: 19, col: 0, oldText: \"31\u00d782=\", newText: \"35\u00d711=\" },\n  { row: 19, col: 1, oldText: \"98\u00d721=\", newText: \"99\u00d797=\" },\n  { row: 19, col: 2, oldText: \"66\u00d782=\", newText: \"14\u00d754=\" },\n  { row: 19, col: 3, oldText: \"55\u00d756=\", newText: \"15\u00d730=\" },\n  { row: 19, col: 4, oldText: \"31\u00d755=\", newText: \"27\u00d784=\" },\n];\n\nfor (const { row, col, oldText, newText } of cellReplacements) {\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchCase: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Cell (${row}, ${col}) text not found: ${oldText}`);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update master to output generated at c8c62b6\n#\n# Updates the date line and each of the 25 multiplication-problem cells in\n# the table. Cells/paragraphs are targeted by POSITION (paragraph index /\n# table row+column), and each Range's trailing paragraph-mark (and, for table\n# cells, the following end-of-cell mark) is trimmed off before the literal\n# text is overwritten. This keeps every edit scoped to exactly one run and\n# preserves its rPr/pPr, and it sidesteps Find/Execute, whose match here is\n# not reliably constrained to the Range it is called on. That matters because\n# several of the new values collide with other cells' original text (e.g. the\n# cell that becomes \"62\u00d731=\" sits beside a different cell that originally\n# held \"62\u00d731=\" and is being changed to something else), so a naive\n# whole-document find/replace (or an under-scoped one) would edit the wrong cell.\n\n$d = $word.ActiveDocument\n\nfunction Set-RangeText($range, $newText) {\n  $r = $range.Duplicate\n  $r.MoveEnd(1, -1) | Out-Null   # wdCharacter; drops the trailing (paragraph/cell) mark(s)\n  $r.Text = $newText\n}\n\nSet-RangeText $d.Paragraphs.Item(1).Range \"2025-10-30 Thursday\"\n\n$table = $d.Tables.Item(1)\n\n$cellReplacements = @(\n  @{ Row = 1; Col = 1; NewText = \"72\u00d735=\" },  # was \"49\u00d764=\"\n  @{ Row = 1; Col = 2; NewText = \"52\u00d780=\" },  # was \"23\u00d738=\"\n  @{ Row = 1; Col = 3; NewText = \"62\u00d731=\" },  # was \"40\u00d797=\"\n  @{ Row = 1; Col = 4; NewText = \"84\u00d736=\" },  # was \"58\u00d786=\"\n  @{ Row = 1; Col = 5; NewText = \"99\u00d772=\" },  # was \"75\u00d736=\"\n  @{ Row = 5; Col = 1; NewText = \"28\u00d738=\" },  # was \"64\u00d747=\"\n  @{ Row = 5; Col = 2; NewText = \"19\u00d745=\" },  # was \"18\u00d714=\"\n  @{ Row = 5; Col = 3; NewText = \"48\u00d755=\" },  # was \"89\u00d722=\"\n  @{ Row = 5; Col = 4; NewText = \"80\u00d778=\" },  # was \"54\u00d743=\"\n  @{ Row = 5; Col = 5; NewText = \"89\u00d739=\" },  # was \"89\u00d745=\"\n  @{ Row = 10; Col = 1; NewText = \"31\u00d781=\" },  # was \"31\u00d724=\"\n  @{ Row = 10; Col = 2; NewText = \"88\u00d799=\" },  # was \"33\u00d787=\"\n  @{ Row = 10; Col = 3; NewText = \"27\u00d735=\" },  # was \"71\u00d781=\"\n  @{ Row = 10; Col = 4; NewText = \"26\u00d767=\" },  # was \"36\u00d799=\"\n  @{ Row = 10; Col = 5; NewText = \"28\u00d719=\" },  # was \"49\u00d770=\"\n  @{ Row = 15; Col = 1; NewText = \"92\u00d798=\" },  # was \"32\u00d720=\"\n  @{ Row = 15; Col = 2; NewText = \"34\u00d771=\" },  # was \"62\u00d731=\"\n  @{ Row = 15; Col = 3; NewText = \"18\u00d797=\" },  # was \"63\u00d762=\"\n  @{ Row = 15; Col = 4; NewText = \"80\u00d724=\" },  # was \"12\u00d790=\"\n  @{ Row = 15; Col = 5; NewText = \"86\u00d791=\" },  # was \"25\u00d773=\"\n  @{ Row = 20; Col = 1; NewText = \"35\u00d711=\" },  # was \"31\u00d782=\"\n  @{ Row = 20; Col = 2; NewText = \"99\u00d797=\" },  # was \"98\u00d721=\"\n  @{ Row = 20; Col = 3; NewText = \"14\u00d754=\" },  # was \"66\u00d782=\"\n  @{ Row = 20; Col = 4; NewText = \"15\u00d730=\" },  # was \"55\u00d756=\"\n  @{ Row = 20; Col = 5; NewText = \"27\u00d784=\" },  # was \"31\u00d755=\"\n)\n\nforeach ($item in $cellReplacements) {\n  $cell = $table.Cell($item.Row, $item.Col)\n  Set-RangeText $cell.Range $item.NewText\n}\n"}
